# Fix a test value: the CITY cell for the row with Kelly Smithson (row 3)
# was "Ottawa" but should be "Gloucester".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = "Gloucester"

# Update the selected/active cell to match the saved file state.
$ws.Range("H2").Select()
